$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 337.14285
$ws.Range("I2").Value = 433.7143
$ws.Range("K2").Value = 433.7143
$ws.Range("M2").Value = -320.7143
$ws.Range("H18").Value = 640
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568
$ws.Range("H58").Value = 7707.8
$ws.Range("I58").Value = 513.3333
$ws.Range("K58").Value = 1539.9999
$ws.Range("M58").Value = -1389.9999
$ws.Range("H74").Value = 5078.231
$ws.Range("I74").Value = 3913.4
$ws.Range("K74").Value = 3913.4
$ws.Range("M74").Value = -2977.4
$ws.Range("H77").Value = 5078.231
$ws.Range("I77").Value = 3913.4
$ws.Range("K77").Value = 19567
$ws.Range("M77").Value = -14887
$ws.Range("H101").Value = 774.6923
$ws.Range("I101").Value = 905.875
$ws.Range("J101").Value = 564.8
$ws.Range("K101").Value = 2717.625
$ws.Range("L101").Value = 1694.4
$ws.Range("M101").Value = -1095.625
$ws.Range("N101").Value = -4938.4
$ws.Range("H106").Value = 3000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 3000
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -4262
$ws.Range("H135").Value = 5055.4
$ws.Range("I135").Value = 5749.5835
$ws.Range("J135").Value = 2278.6667
$ws.Range("K135").Value = 51746.2515
$ws.Range("L135").Value = 20508.0003
$ws.Range("M135").Value = -49211.2515
$ws.Range("N135").Value = -25578.0003

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 500024500
$ws.Range("I43").Value = 1000000000
$ws.Range("K43").Value = 1000000000
$ws.Range("M43").Value = -999999687
$ws.Range("H45").Value = 1686
$ws.Range("I45").Value = 1648
$ws.Range("K45").Value = 1648
$ws.Range("M45").Value = -1271
$ws.Range("H61").Value = 3333.4546
$ws.Range("I61").Value = 2351.1875
$ws.Range("J61").Value = 5952.8335
$ws.Range("K61").Value = 2351.1875
$ws.Range("L61").Value = 5952.8335
$ws.Range("M61").Value = -2139.1875
$ws.Range("N61").Value = -6376.8335
$ws.Range("H63").Value = 2226
$ws.Range("I63").Value = 2226
$ws.Range("K63").Value = 2226
$ws.Range("M63").Value = -1540
$ws.Range("H66").Value = 2226
$ws.Range("I66").Value = 2226
$ws.Range("K66").Value = 11130
$ws.Range("M66").Value = -7698
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 27188.91
$ws.Range("I74").Value = 31812.264
$ws.Range("J74").Value = 2090.7144
$ws.Range("K74").Value = 31812.264
$ws.Range("L74").Value = 2090.7144
$ws.Range("M74").Value = -30938.264
$ws.Range("N74").Value = -3838.7144
$ws.Range("H77").Value = 27188.91
$ws.Range("I77").Value = 31812.264
$ws.Range("J77").Value = 2090.7144
$ws.Range("K77").Value = 159061.32
$ws.Range("L77").Value = 10453.572
$ws.Range("M77").Value = -154693.32
$ws.Range("N77").Value = -19189.572
$ws.Range("H132").Value = 55397.445
$ws.Range("I132").Value = 1595.4117
$ws.Range("J132").Value = 221694.64
$ws.Range("K132").Value = 4786.2351
$ws.Range("L132").Value = 665083.92
$ws.Range("M132").Value = -2256.2351
$ws.Range("N132").Value = -670143.92
$ws.Range("H133").Value = 117500
$ws.Range("J133").Value = 117500
$ws.Range("L133").Value = 117500
$ws.Range("N133").Value = -122560
$ws.Range("H134").Value = 97886.25
$ws.Range("J134").Value = 97886.25
$ws.Range("L134").Value = 97886.25
$ws.Range("N134").Value = -108026.25
$ws.Range("H136").Value = 3333.4546
$ws.Range("I136").Value = 2351.1875
$ws.Range("J136").Value = 5952.8335
$ws.Range("K136").Value = 7053.5625
$ws.Range("L136").Value = 17858.5005
$ws.Range("M136").Value = -4503.5625
$ws.Range("N136").Value = -22958.5005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3032.2415
$ws.Range("I134").Value = 2343.8462
$ws.Range("K134").Value = 7031.5386
$ws.Range("M134").Value = -4496.5386

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3249.2856
$ws.Range("I132").Value = 3268.4614
$ws.Range("K132").Value = 9805.3842
$ws.Range("M132").Value = -7275.3842

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 376.5
$ws.Range("J11").Value = 1999
$ws.Range("L11").Value = 5997
$ws.Range("N11").Value = -6277
$ws.Range("H16").Value = 287.66666
$ws.Range("J16").Value = 376
$ws.Range("L16").Value = 1128
$ws.Range("N16").Value = -1474
$ws.Range("H19").Value = 750
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 1500
$ws.Range("M19").Value = -1326
$ws.Range("H92").Value = 988
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 982
$ws.Range("K92").Value = 3000
$ws.Range("L92").Value = 2946
$ws.Range("M92").Value = -1752
$ws.Range("N92").Value = -5442
$ws.Range("H138").Value = 2289.8572
$ws.Range("I138").Value = 2005.8
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 6017.4
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -877.3999999999996
$ws.Range("N138").Value = -19280

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10910.2
$ws.Range("I99").Value = 10910.2
$ws.Range("K99").Value = 10910.2
$ws.Range("M99").Value = -8664.200000000001
$ws.Range("H113").Value = 3503.3914
$ws.Range("I113").Value = 2776.6428
$ws.Range("J113").Value = 4633.8887
$ws.Range("K113").Value = 2776.6428
$ws.Range("L113").Value = 4633.8887
$ws.Range("M113").Value = -606.6428000000001
$ws.Range("N113").Value = -8973.8887

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 12000
$ws.Range("I20").Value = 12000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -11774
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 2380.2173
$ws.Range("J46").Value = 5343.125
$ws.Range("L46").Value = 5343.125
$ws.Range("N46").Value = -5719.125
$ws.Range("H55").Value = 850.38464
$ws.Range("I55").Value = 388
$ws.Range("K55").Value = 388
$ws.Range("M55").Value = -215
$ws.Range("H132").Value = 8079.524
$ws.Range("I132").Value = 4119.3335
$ws.Range("J132").Value = 17980
$ws.Range("K132").Value = 12358.0005
$ws.Range("L132").Value = 53940
$ws.Range("M132").Value = -9828.000499999998
$ws.Range("N132").Value = -59000
$ws.Range("H136").Value = 3050.3667
$ws.Range("I136").Value = 2944.1304
$ws.Range("K136").Value = 8832.3912
$ws.Range("M136").Value = -6282.3912

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2749.75
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 2749.75
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 73625.28999999999
$ws.Range("I122").Value = 144057.28
$ws.Range("K122").Value = 432171.84
$ws.Range("M122").Value = -429721.84
$ws.Range("H132").Value = 2463.6924
$ws.Range("I132").Value = 2161.9092
$ws.Range("J132").Value = 4123.5
$ws.Range("K132").Value = 6485.7276
$ws.Range("L132").Value = 12370.5
$ws.Range("M132").Value = -3955.7276
$ws.Range("N132").Value = -17430.5
$ws.Range("H135").Value = 12588379
$ws.Range("J135").Value = 12588379
$ws.Range("L135").Value = 12588379
$ws.Range("N135").Value = -12598519
$ws.Range("H136").Value = 12330.9375
$ws.Range("I136").Value = 11460
$ws.Range("K136").Value = 34380
$ws.Range("M136").Value = -31830
